$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Existing row 13 ("CLOck") is updated to use the new "CLOCK" text.
$ws.Range("F13").Value = "CLOCK"

# New row 65: scroll wheel / channel object text entry.
$ws.Range("B65").Value = "SingleUseId65"
$ws.Range("C65").Value = "Default"
$ws.Range("D65").Value = "Center"
$ws.Range("E65").Value = "LTR"
$ws.Range("F65").Value = "<value>"
